# Fruta / hortaliza, semanal
# Inserts two new weekly price rows for "Espárragos" (Mapocho Venta Directa
# de Santiago) at rows 29-30, pushing the existing rows 29-47 down to 31-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 29, shifting the rest of the
# table (previously rows 29:47) down to rows 31:49.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(30).Insert()

# New row 29: Banquete quality, week of 2021-10-19
$ws.Range("A29").Value = 12
$ws.Range("B29").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44488
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 300000000
$ws.Range("G29").Value = "Espárragos"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Banquete"
$ws.Range("J29").Value = 370
$ws.Range("K29").Value = 1200
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = 1200
$ws.Range("N29").Value = "$/kilo"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 1200
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"

# New row 30: Primera quality, week of 2021-10-19
$ws.Range("A30").Value = 12
$ws.Range("B30").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 44488
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = 300000000
$ws.Range("G30").Value = "Espárragos"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 480
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = 1000
$ws.Range("N30").Value = "$/kilo"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1000
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"
